# Add a "player_team" column (F) and assign each player to the "PPL" team.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the neighbouring "player_position" column (E) onto
# the new column (F) so the new header/value cells share the same style as
# the rest of the header/data row instead of picking up a brand new style.
$ws.Range("E1:E3").Copy()
$ws.Range("F1:F3").PasteSpecial(-4122)

# New header + data.
$ws.Range("F1").Value = "player_team"
$ws.Range("F2").Value = "PPL"
$ws.Range("F3").Value = "PPL"
